$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Approved/Rejected column (I) and ReasonToReject column (J)
# for rows 2, 8 and 10: change "Approved" -> "Rejected" and set ReasonToReject -> "Nil"
$ws.Range("I2").Value = "Rejected"
$ws.Range("J2").Value = "Nil"

$ws.Range("I8").Value = "Rejected"
$ws.Range("J8").Value = "Nil"

$ws.Range("I10").Value = "Rejected"
$ws.Range("J10").Value = "Nil"

# Update the current selection to match the saved view state (J10)
$ws.Range("J10").Select()
